$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 66
$ws.Range("A66").Value = 44700
$ws.Range("B66").Value = 0.33333333333333331
$ws.Range("C66").Value = 0.39305555555555555
$ws.Range("E66").Value = "Changement de stratégie"
$ws.Range("F66").Value = "L'algorithme est très différent dans sont fonctionnement que les algorithme vu précédemment`n"
$ws.Range("G66").Value = "Il faudra faire l'algorithme pas à pas ajoutant les complexité une a une afin de:`n1) etre plus efficace dans son implémnentation`n2) proposer une documentation plus structuré et améliorer la compréhension"

# Row 67
$ws.Range("A67").Value = 44700
$ws.Range("B67").Value = 0.39305555555555555
$ws.Range("C67").Value = 0.51041666666666663
$ws.Range("E67").Value = "HPA: Conception D'un chunk"
$ws.Range("F67").Value = "But:`nMontrer la logique des emplacement des portes et les chemins qui les relient entre elle sans se préoccuper des partitions adjacentes pour le moment"

# Row 68
$ws.Range("A68").Value = 44700
$ws.Range("B68").Value = 0.5625
$ws.Range("C68").Value = 0.65416666666666667
$ws.Range("E68").Value = "Analyse et documentation des besoins pour l'algorithme"

# Row 69
$ws.Range("A69").Value = 44700
$ws.Range("B69").Value = 0.65416666666666667
$ws.Range("C69").Value = 0.70486111111111116
$ws.Range("E69").Value = "implémentations"
$ws.Range("F69").Value = "Quelque soucis dans les container utilisé"

# Update sheet view to reflect scroll/selection position (journal grew, so the
# author had scrolled further down and was editing near F70 when they saved)
$ws.Application.ActiveWindow.ScrollRow = 64
$ws.Range("F70").Select()

